$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 11.629057
$ws.Range("H2").Value = 34.887171
$ws.Range("I2").Value = 0.3062678464977661
$ws.Range("J2").Value = 0.3062678464977662
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 10.41202966666667
$ws.Range("N2").Value = 31.236089
$ws.Range("Q2").Value = 121.0820864793577
$ws.Range("R2").Value = 1089.738778314219
$ws.Range("S2").Value = 0.3062678464977661
$ws.Range("T2").Value = 0.3062678464977662

# Row 3
$ws.Range("I3").Value = 0.6269156120645606
$ws.Range("J3").Value = 0.6269156120645607
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 10.41202966666667
$ws.Range("N3").Value = 31.236089
$ws.Range("Q3").Value = 247.8492313943058
$ws.Range("R3").Value = 2230.643082548752
$ws.Range("S3").Value = 0.6269156120645606
$ws.Range("T3").Value = 0.6269156120645607

# Row 4
$ws.Range("G4").Value = 2.537038666666667
$ws.Range("H4").Value = 7.611116
$ws.Range("I4").Value = 0.06681654143767324
$ws.Range("J4").Value = 0.06681654143767324
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 10.41202966666667
$ws.Range("N4").Value = 31.236089
$ws.Range("Q4").Value = 26.41572186281378
$ws.Range("R4").Value = 237.741496765324
$ws.Range("S4").Value = 0.06681654143767324
$ws.Range("T4").Value = 0.06681654143767324
